$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.995.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.46'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6252'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9988'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07573'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2918'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.50'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.834.33'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.945'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6624'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001019'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +18.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.037'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.997.26'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '225.92'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.33'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9992'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.479'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1376'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.87'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.488'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.089'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.006'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.184'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05231'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7332'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.685'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.235.16'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.757'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01779'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.316'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8971'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9985'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.77'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.978.03'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000125'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.09'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5102'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4030'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.842'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05740'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.655'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.24%  '
